# Bugs and improvements sheet update
# - Replace the placeholder "Video/Image Link" values (SS / V) with real
#   Google-Drive hyperlinks in column E
# - Add a second "Steps to Recreate" bullet (D3) for the new discount-format row
# - Resize things (row 3 height, column E width) to fit the new content
# - Remove the two trailing blank rows
# - Update the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Steps to Recreate" text for row 3 (Non uniform Discount format) ---
$ws.Range("D3").Value = "1. Go to any category or section`n2. Some products have Discount price while there's discount percentage in others"
$ws.Range("D3").WrapText = $true

# --- Video/Image Link column: real hyperlinked URLs instead of SS/V placeholders ---
$ws.Range("E2").Value = "https://drive.google.com/file/d/1szzErH2C_SeYx8dBQnBm17Fh6eUC4X-V/view?usp=drive_link"
$ws.Range("E2").VerticalAlignment = -4108  # xlCenter
$ws.Hyperlinks.Add($ws.Range("E2"), "https://drive.google.com/file/d/1szzErH2C_SeYx8dBQnBm17Fh6eUC4X-V/view?usp=drive_link")

$ws.Range("E3").Value = "https://drive.google.com/file/d/1CqYrWFFT0XQ8tdBdk14XpzIEf1l7xKaf/view?usp=drive_link"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://drive.google.com/file/d/1CqYrWFFT0XQ8tdBdk14XpzIEf1l7xKaf/view?usp=drive_link")

$ws.Range("E4").Value = "https://drive.google.com/file/d/1QleOHP9UNqrotyynYyLmnJcccMPow_R7/view?usp=drive_link"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://drive.google.com/file/d/1QleOHP9UNqrotyynYyLmnJcccMPow_R7/view?usp=drive_link")

$ws.Range("E5").Value = "https://drive.google.com/file/d/12otaMISfL_ZA1ZJZAfrgk_m7H2egPVGe/view?usp=drive_link"
$ws.Range("E5").VerticalAlignment = -4108  # xlCenter
$ws.Hyperlinks.Add($ws.Range("E5"), "https://drive.google.com/file/d/12otaMISfL_ZA1ZJZAfrgk_m7H2egPVGe/view?usp=drive_link")

# --- Row height / column width tweaks to fit the new wrapped text / links ---
$ws.Rows(3).RowHeight = 47.25
$ws.Columns(5).ColumnWidth = 29.75

# --- Drop the two trailing empty rows now that the sheet ends at row 5 ---
$ws.Rows("6:7").Delete()

# --- Update current selection / view ---
$ws.Range("D17").Select()
